$d = $word.ActiveDocument

# Trim the trailing clause from the ASICs/mining-centralization sentence:
# "...no ASICs will be created for quite some time, until Smartcash reaches
#  a considerable market cap." -> "...no ASICs will be created for quite some time."
$found = $d.Content.Find.Execute(
    ", until Smartcash reaches a considerable market cap.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ".", 2)

if (-not $found) {
    throw "Could not locate the sentence to trim."
}
